$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the timestamp in A18 (minor precision fix from the data source)
$ws.Range("A18").Value = 44331.77618375231

# Append the new day's row of data (row 19)
$ws.Range("A19").Value = 44332.77772781081
$ws.Range("B19").Value = 73794
$ws.Range("C19").Value = 62153
$ws.Range("D19").Value = 3196
$ws.Range("E19").Value = 2102
$ws.Range("F19").Value = 1490
$ws.Range("G19").Value = 19195
$ws.Range("H19").Value = 1301
$ws.Range("I19").Value = 868
$ws.Range("J19").Value = 204
